# Re-orders the "record" (species observation) data held in a handful of
# rows of the Artfynd sheet. The rows themselves stay put (same Id-column
# position in the sheet layout, same shared location/metadata columns such
# as Lokalnamn, Ost/Nord-adjacent county data, dates, observer, etc.) but
# the species-specific payload -- Id, Taxonsorteringsordning, TaxonId,
# Artnamn, Vetenskapligt namn, Auktor, Enhet, Kön, Aktivitet, Ost, Nord,
# Publik kommentar, Bestämningsmetod, Substrat, Substrat-beskrivning --
# is relocated between rows 15<->16, 19<->20, 21<->22 and rotated through
# 24->25->26->24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that carry the record-specific payload that needs to move.
$cols = @("A","B","E","F","G","H","J","L","M","Q","R","AC","AF","AM","AO")

function Get-RowData($row) {
    $data = @{}
    foreach ($c in $cols) {
        $data[$c] = $ws.Range("$c$row").Value()
    }
    return $data
}

function Set-RowData($row, $data) {
    foreach ($c in $cols) {
        $v = $data[$c]
        if ($null -eq $v) { $v = "" }
        $ws.Range("$c$row").Value = $v
    }
}

# Snapshot every source row BEFORE any writes happen. This matters for the
# three-way rotation among rows 24/25/26 where each row is both a source
# and a destination.
$row15 = Get-RowData 15
$row16 = Get-RowData 16
$row19 = Get-RowData 19
$row20 = Get-RowData 20
$row21 = Get-RowData 21
$row22 = Get-RowData 22
$row24 = Get-RowData 24
$row25 = Get-RowData 25
$row26 = Get-RowData 26

# Swap pairs.
Set-RowData 15 $row16
Set-RowData 16 $row15

Set-RowData 19 $row20
Set-RowData 20 $row19

Set-RowData 21 $row22
Set-RowData 22 $row21

# Three-way rotation: 24 receives 26's data, 25 receives 24's (original)
# data, 26 receives 25's (original) data.
Set-RowData 24 $row26
Set-RowData 25 $row24
Set-RowData 26 $row25
